$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 249; this shifts the existing rows 249-303 down to 250-304
$ws.Rows(249).Insert()

# Populate the newly inserted row 249 with the new data point
$ws.Range("A249").Value2 = 5
$ws.Range("B249").Value2 = "Macroferia Regional de Talca"
$ws.Range("C249").Value2 = "Maule"
$ws.Range("D249").Value2 = 44641
$ws.Range("D249").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E249").Value2 = 7
$ws.Range("F249").Value2 = 100114013
$ws.Range("G249").Value2 = "Zanahoria"
$ws.Range("H249").Value2 = "Sin especificar"
$ws.Range("I249").Value2 = "Primera"
$ws.Range("J249").Value2 = 500
$ws.Range("K249").Value2 = 7000
$ws.Range("L249").Value2 = 7000
$ws.Range("M249").Value2 = 7000
$ws.Range("N249").Value2 = "$/saco 20 kilos"
$ws.Range("O249").Value2 = "Región de Ñuble"
$ws.Range("P249").Value2 = 350
$ws.Range("Q249").Value2 = 20
$ws.Range("R249").Value2 = "Hortaliza"
